# Update the CDCF workbook from the "EU model" variant to the "US model" variant.
$wb = $excel.ActiveWorkbook

# --- Sheet "About" ---
$about = $wb.Worksheets.Item("About")

# Replace the "For the EU model..." block with the US-model text, and drop the
# two trailing informational rows (ton-mile conversion note + helper formula).
$about.Range("A11").Value = "For the U.S. model, the desired output units are:"
$about.Range("A12").Value = "trillion passenger-miles"
$about.Range("A13").Value = "trillion freight ton-miles"
$about.Range("A15").ClearContents()
$about.Range("A17").Clear()

# --- Sheet "CDCF-PMpPDOU" ---
$pm = $wb.Worksheets.Item("CDCF-PMpPDOU")
$pm.Range("B2").Formula = "=10^12"

# --- Sheet "CDCF-FTMpFDOU" ---
$ftm = $wb.Worksheets.Item("CDCF-FTMpFDOU")
$ftm.Range("B2").Formula = "=10^12"
$ftm.Range("B5").Clear()
